$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Put the Neo4j/Cypher query used to generate this test-case sheet into A2
# (new shared string), keep it on the "Normal 2" (wrap-text) style already
# applied to A2, and size the row to fit the multi-line query text.
$ws.Range("A2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Adenocarcinoma, NOS''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'
$ws.Range("A2").RowHeight = 87

# Re-point the selection at the data that was just (re)generated.
$ws.Range("B2:B6").Select()
